$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Cells.Item(27, 2).Value2 = 2594729
$ws.Cells.Item(27, 6).Value2 = "FC Nordsjaelland"
$ws.Cells.Item(27, 7).Value2 = "AaB"
$ws.Cells.Item(27, 10).Value2 = "D"
$ws.Cells.Item(27, 9).Value2 = 2
$ws.Cells.Item(27, 11).Value2 = 2.45
$ws.Cells.Item(27, 12).Value2 = 3.4
$ws.Cells.Item(27, 13).Value2 = 2.75
$ws.Cells.Item(27, 14).Value2 = 2.9
$ws.Cells.Item(27, 15).Value2 = 3.4
$ws.Cells.Item(27, 16).Value2 = 2.25
$ws.Cells.Item(27, 17).Value2 = 0.25
$ws.Cells.Item(27, 18).Value2 = 1.85
$ws.Cells.Item(27, 19).Value2 = 2
$ws.Cells.Item(27, 20).Value2 = 2.75
$ws.Cells.Item(27, 21).Value2 = 1.925
$ws.Cells.Item(27, 22).Value2 = 1.925
$ws.Cells.Item(27, 23).Value2 = -1
$ws.Cells.Item(27, 24).Value2 = 2.4
$ws.Cells.Item(27, 26).Value2 = 0.425
$ws.Cells.Item(27, 27).Value2 = -0.5
$ws.Cells.Item(27, 28).Value2 = 0.925

# Row 28
$ws.Cells.Item(28, 2).Value2 = 2594756
$ws.Cells.Item(28, 6).Value2 = "Odense BK"
$ws.Cells.Item(28, 7).Value2 = "Randers FC"
$ws.Cells.Item(28, 10).Value2 = "H"
$ws.Cells.Item(28, 9).Value2 = 1
$ws.Cells.Item(28, 11).Value2 = 2.5
$ws.Cells.Item(28, 12).Value2 = 3.3
$ws.Cells.Item(28, 13).Value2 = 2.8
$ws.Cells.Item(28, 14).Value2 = 2.8
$ws.Cells.Item(28, 15).Value2 = 3.3
$ws.Cells.Item(28, 16).Value2 = 2.45
$ws.Cells.Item(28, 17).Value2 = 0
$ws.Cells.Item(28, 18).Value2 = 2
$ws.Cells.Item(28, 19).Value2 = 1.85
$ws.Cells.Item(28, 20).Value2 = 2.5
$ws.Cells.Item(28, 21).Value2 = 2.025
$ws.Cells.Item(28, 22).Value2 = 1.825
$ws.Cells.Item(28, 23).Value2 = 1.8
$ws.Cells.Item(28, 24).Value2 = -1
$ws.Cells.Item(28, 26).Value2 = 1
$ws.Cells.Item(28, 27).Value2 = -1
$ws.Cells.Item(28, 28).Value2 = 1.025

# Row 32
$ws.Cells.Item(32, 2).Value2 = 2594732
$ws.Cells.Item(32, 6).Value2 = "AGF Aarhus"
$ws.Cells.Item(32, 7).Value2 = "FC Nordsjaelland"
$ws.Cells.Item(32, 10).Value2 = "A"
$ws.Cells.Item(32, 8).Value2 = 0
$ws.Cells.Item(32, 11).Value2 = 1.666
$ws.Cells.Item(32, 12).Value2 = 4
$ws.Cells.Item(32, 13).Value2 = 4.75
$ws.Cells.Item(32, 14).Value2 = 1.533
$ws.Cells.Item(32, 15).Value2 = 4
$ws.Cells.Item(32, 16).Value2 = 5.5
$ws.Cells.Item(32, 17).Value2 = -1
$ws.Cells.Item(32, 18).Value2 = 1.95
$ws.Cells.Item(32, 19).Value2 = 1.9
$ws.Cells.Item(32, 20).Value2 = 2.75
$ws.Cells.Item(32, 21).Value2 = 1.85
$ws.Cells.Item(32, 22).Value2 = 2
$ws.Cells.Item(32, 23).Value2 = -1
$ws.Cells.Item(32, 25).Value2 = 4.5
$ws.Cells.Item(32, 26).Value2 = -1
$ws.Cells.Item(32, 27).Value2 = 0.8999999999999999
$ws.Cells.Item(32, 28).Value2 = -1
$ws.Cells.Item(32, 29).Value2 = 1

# Row 33
$ws.Cells.Item(33, 2).Value2 = 2594734
$ws.Cells.Item(33, 6).Value2 = "AC Horsens"
$ws.Cells.Item(33, 7).Value2 = "AaB"
$ws.Cells.Item(33, 10).Value2 = "H"
$ws.Cells.Item(33, 8).Value2 = 2
$ws.Cells.Item(33, 11).Value2 = 3.2
$ws.Cells.Item(33, 12).Value2 = 3.5
$ws.Cells.Item(33, 13).Value2 = 2.2
$ws.Cells.Item(33, 14).Value2 = 3.25
$ws.Cells.Item(33, 15).Value2 = 3.5
$ws.Cells.Item(33, 16).Value2 = 2.2
$ws.Cells.Item(33, 17).Value2 = 0.25
$ws.Cells.Item(33, 18).Value2 = 1.925
$ws.Cells.Item(33, 19).Value2 = 1.925
$ws.Cells.Item(33, 20).Value2 = 2.25
$ws.Cells.Item(33, 21).Value2 = 1.95
$ws.Cells.Item(33, 22).Value2 = 1.9
$ws.Cells.Item(33, 23).Value2 = 2.25
$ws.Cells.Item(33, 25).Value2 = -1
$ws.Cells.Item(33, 26).Value2 = 0.925
$ws.Cells.Item(33, 27).Value2 = -1
$ws.Cells.Item(33, 28).Value2 = 0.95
$ws.Cells.Item(33, 29).Value2 = -1

# Row 80
$ws.Cells.Item(80, 2).Value2 = 3377638
$ws.Cells.Item(80, 6).Value2 = "AC Horsens"
$ws.Cells.Item(80, 7).Value2 = "AaB"
$ws.Cells.Item(80, 11).Value2 = 3.5
$ws.Cells.Item(80, 13).Value2 = 2.05
$ws.Cells.Item(80, 14).Value2 = 3.9
$ws.Cells.Item(80, 15).Value2 = 3.5
$ws.Cells.Item(80, 16).Value2 = 1.95
$ws.Cells.Item(80, 17).Value2 = 0.5
$ws.Cells.Item(80, 18).Value2 = 1.875
$ws.Cells.Item(80, 19).Value2 = 1.975
$ws.Cells.Item(80, 20).Value2 = 2.5
$ws.Cells.Item(80, 21).Value2 = 1.9
$ws.Cells.Item(80, 22).Value2 = 1.95
$ws.Cells.Item(80, 23).Value2 = 2.9
$ws.Cells.Item(80, 26).Value2 = 0.875
$ws.Cells.Item(80, 29).Value2 = 0.95

# Row 81
$ws.Cells.Item(81, 2).Value2 = 3377637
$ws.Cells.Item(81, 6).Value2 = "Sonderjyske"
$ws.Cells.Item(81, 7).Value2 = "Vejle"
$ws.Cells.Item(81, 11).Value2 = 2.35
$ws.Cells.Item(81, 13).Value2 = 2.9
$ws.Cells.Item(81, 14).Value2 = 2.4
$ws.Cells.Item(81, 15).Value2 = 2.9
$ws.Cells.Item(81, 16).Value2 = 3.3
$ws.Cells.Item(81, 17).Value2 = -0.25
$ws.Cells.Item(81, 18).Value2 = 2.05
$ws.Cells.Item(81, 19).Value2 = 1.8
$ws.Cells.Item(81, 20).Value2 = 2.25
$ws.Cells.Item(81, 21).Value2 = 2.05
$ws.Cells.Item(81, 22).Value2 = 1.8
$ws.Cells.Item(81, 23).Value2 = 1.4
$ws.Cells.Item(81, 26).Value2 = 1.05
$ws.Cells.Item(81, 29).Value2 = 0.8

# Row 110
$ws.Cells.Item(110, 2).Value2 = 3377653
$ws.Cells.Item(110, 6).Value2 = "Sonderjyske"
$ws.Cells.Item(110, 7).Value2 = "AaB"
$ws.Cells.Item(110, 10).Value2 = "A"
$ws.Cells.Item(110, 8).Value2 = 0
$ws.Cells.Item(110, 9).Value2 = 4
$ws.Cells.Item(110, 11).Value2 = 2.7
$ws.Cells.Item(110, 12).Value2 = 3.6
$ws.Cells.Item(110, 13).Value2 = 2.35
$ws.Cells.Item(110, 14).Value2 = 3.3
$ws.Cells.Item(110, 15).Value2 = 3.5
$ws.Cells.Item(110, 16).Value2 = 2.1
$ws.Cells.Item(110, 17).Value2 = 0.25
$ws.Cells.Item(110, 18).Value2 = 1.95
$ws.Cells.Item(110, 19).Value2 = 1.9
$ws.Cells.Item(110, 21).Value2 = 2
$ws.Cells.Item(110, 22).Value2 = 1.85
$ws.Cells.Item(110, 23).Value2 = -1
$ws.Cells.Item(110, 25).Value2 = 1.1
$ws.Cells.Item(110, 26).Value2 = -1
$ws.Cells.Item(110, 27).Value2 = 0.8999999999999999
$ws.Cells.Item(110, 28).Value2 = 1

# Row 112
$ws.Cells.Item(112, 2).Value2 = 3377654
$ws.Cells.Item(112, 6).Value2 = "Odense BK"
$ws.Cells.Item(112, 7).Value2 = "AC Horsens"
$ws.Cells.Item(112, 10).Value2 = "H"
$ws.Cells.Item(112, 8).Value2 = 4
$ws.Cells.Item(112, 9).Value2 = 0
$ws.Cells.Item(112, 11).Value2 = 1.5
$ws.Cells.Item(112, 12).Value2 = 4.2
$ws.Cells.Item(112, 13).Value2 = 6
$ws.Cells.Item(112, 14).Value2 = 1.55
$ws.Cells.Item(112, 15).Value2 = 4
$ws.Cells.Item(112, 16).Value2 = 5.75
$ws.Cells.Item(112, 17).Value2 = -1
$ws.Cells.Item(112, 18).Value2 = 2.025
$ws.Cells.Item(112, 19).Value2 = 1.825
$ws.Cells.Item(112, 21).Value2 = 1.9
$ws.Cells.Item(112, 22).Value2 = 1.95
$ws.Cells.Item(112, 23).Value2 = 0.55
$ws.Cells.Item(112, 25).Value2 = -1
$ws.Cells.Item(112, 26).Value2 = 1.025
$ws.Cells.Item(112, 27).Value2 = -1
$ws.Cells.Item(112, 28).Value2 = 0.8999999999999999

# Row 207
$ws.Cells.Item(207, 2).Value2 = 3613141
$ws.Cells.Item(207, 6).Value2 = "Randers FC"
$ws.Cells.Item(207, 7).Value2 = "Vejle"
$ws.Cells.Item(207, 10).Value2 = "H"
$ws.Cells.Item(207, 8).Value2 = 4
$ws.Cells.Item(207, 11).Value2 = 1.55
$ws.Cells.Item(207, 12).Value2 = 4
$ws.Cells.Item(207, 13).Value2 = 5.5
$ws.Cells.Item(207, 14).Value2 = 1.533
$ws.Cells.Item(207, 15).Value2 = 4.2
$ws.Cells.Item(207, 16).Value2 = 6.5
$ws.Cells.Item(207, 17).Value2 = -1
$ws.Cells.Item(207, 18).Value2 = 1.975
$ws.Cells.Item(207, 19).Value2 = 1.875
$ws.Cells.Item(207, 20).Value2 = 2.75
$ws.Cells.Item(207, 21).Value2 = 2
$ws.Cells.Item(207, 22).Value2 = 1.85
$ws.Cells.Item(207, 23).Value2 = 0.5329999999999999
$ws.Cells.Item(207, 24).Value2 = -1
$ws.Cells.Item(207, 26).Value2 = 0.9750000000000001
$ws.Cells.Item(207, 27).Value2 = -1
$ws.Cells.Item(207, 28).Value2 = 1
$ws.Cells.Item(207, 29).Value2 = -1

# Row 208
$ws.Cells.Item(208, 2).Value2 = 3613140
$ws.Cells.Item(208, 6).Value2 = "FC Nordsjaelland"
$ws.Cells.Item(208, 7).Value2 = "Silkeborg IF"
$ws.Cells.Item(208, 10).Value2 = "D"
$ws.Cells.Item(208, 8).Value2 = 1
$ws.Cells.Item(208, 11).Value2 = 2.75
$ws.Cells.Item(208, 12).Value2 = 3.6
$ws.Cells.Item(208, 13).Value2 = 2.35
$ws.Cells.Item(208, 14).Value2 = 2.9
$ws.Cells.Item(208, 15).Value2 = 3.8
$ws.Cells.Item(208, 16).Value2 = 2.25
$ws.Cells.Item(208, 17).Value2 = 0.25
$ws.Cells.Item(208, 18).Value2 = 1.875
$ws.Cells.Item(208, 19).Value2 = 1.975
$ws.Cells.Item(208, 20).Value2 = 3
$ws.Cells.Item(208, 21).Value2 = 1.975
$ws.Cells.Item(208, 22).Value2 = 1.875
$ws.Cells.Item(208, 23).Value2 = -1
$ws.Cells.Item(208, 24).Value2 = 2.8
$ws.Cells.Item(208, 26).Value2 = 0.4375
$ws.Cells.Item(208, 27).Value2 = -0.5
$ws.Cells.Item(208, 28).Value2 = -1
$ws.Cells.Item(208, 29).Value2 = 0.875

# Row 261
$ws.Cells.Item(261, 2).Value2 = 4811307
$ws.Cells.Item(261, 6).Value2 = "AGF Aarhus"
$ws.Cells.Item(261, 7).Value2 = "Viborg"
$ws.Cells.Item(261, 10).Value2 = "A"
$ws.Cells.Item(261, 8).Value2 = 0
$ws.Cells.Item(261, 9).Value2 = 2
$ws.Cells.Item(261, 11).Value2 = 2.6
$ws.Cells.Item(261, 12).Value2 = 3.3
$ws.Cells.Item(261, 13).Value2 = 2.6
$ws.Cells.Item(261, 14).Value2 = 3.4
$ws.Cells.Item(261, 15).Value2 = 3.3
$ws.Cells.Item(261, 16).Value2 = 2.2
$ws.Cells.Item(261, 17).Value2 = 0.25
$ws.Cells.Item(261, 18).Value2 = 2
$ws.Cells.Item(261, 19).Value2 = 1.9
$ws.Cells.Item(261, 20).Value2 = 2.25
$ws.Cells.Item(261, 21).Value2 = 1.9
$ws.Cells.Item(261, 22).Value2 = 1.95
$ws.Cells.Item(261, 23).Value2 = -1
$ws.Cells.Item(261, 25).Value2 = 1.2
$ws.Cells.Item(261, 26).Value2 = -1
$ws.Cells.Item(261, 27).Value2 = 0.8999999999999999
$ws.Cells.Item(261, 28).Value2 = -0.5
$ws.Cells.Item(261, 29).Value2 = 0.475

# Row 262
$ws.Cells.Item(262, 2).Value2 = 4811306
$ws.Cells.Item(262, 6).Value2 = "Odense BK"
$ws.Cells.Item(262, 7).Value2 = "Vejle"
$ws.Cells.Item(262, 10).Value2 = "H"
$ws.Cells.Item(262, 8).Value2 = 2
$ws.Cells.Item(262, 9).Value2 = 1
$ws.Cells.Item(262, 11).Value2 = 1.769
$ws.Cells.Item(262, 12).Value2 = 3.7
$ws.Cells.Item(262, 13).Value2 = 4.2
$ws.Cells.Item(262, 14).Value2 = 1.65
$ws.Cells.Item(262, 15).Value2 = 4
$ws.Cells.Item(262, 16).Value2 = 5
$ws.Cells.Item(262, 17).Value2 = -0.75
$ws.Cells.Item(262, 18).Value2 = 1.825
$ws.Cells.Item(262, 19).Value2 = 2.025
$ws.Cells.Item(262, 20).Value2 = 2.5
$ws.Cells.Item(262, 21).Value2 = 1.85
$ws.Cells.Item(262, 22).Value2 = 2
$ws.Cells.Item(262, 23).Value2 = 0.6499999999999999
$ws.Cells.Item(262, 25).Value2 = -1
$ws.Cells.Item(262, 26).Value2 = 0.4125
$ws.Cells.Item(262, 27).Value2 = -0.5
$ws.Cells.Item(262, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(262, 29).Value2 = -1

# Row 305
$ws.Cells.Item(305, 2).Value2 = 4811775
$ws.Cells.Item(305, 6).Value2 = "FC Copenhagen"
$ws.Cells.Item(305, 7).Value2 = "AaB"
$ws.Cells.Item(305, 8).Value2 = 3
$ws.Cells.Item(305, 9).Value2 = 0
$ws.Cells.Item(305, 11).Value2 = 1.416
$ws.Cells.Item(305, 12).Value2 = 4.5
$ws.Cells.Item(305, 13).Value2 = 8
$ws.Cells.Item(305, 14).Value2 = 1.444
$ws.Cells.Item(305, 15).Value2 = 5
$ws.Cells.Item(305, 16).Value2 = 6.5
$ws.Cells.Item(305, 17).Value2 = -1.5
$ws.Cells.Item(305, 18).Value2 = 2
$ws.Cells.Item(305, 19).Value2 = 1.8
$ws.Cells.Item(305, 21).Value2 = 1.875
$ws.Cells.Item(305, 22).Value2 = 1.975
$ws.Cells.Item(305, 23).Value2 = 0.444
$ws.Cells.Item(305, 26).Value2 = 1
$ws.Cells.Item(305, 27).Value2 = -1

# Row 306
$ws.Cells.Item(306, 2).Value2 = 4811403
$ws.Cells.Item(306, 6).Value2 = "Brondby"
$ws.Cells.Item(306, 7).Value2 = "Silkeborg IF"
$ws.Cells.Item(306, 8).Value2 = 2
$ws.Cells.Item(306, 9).Value2 = 1
$ws.Cells.Item(306, 11).Value2 = 2.6
$ws.Cells.Item(306, 12).Value2 = 3.6
$ws.Cells.Item(306, 13).Value2 = 2.6
$ws.Cells.Item(306, 14).Value2 = 1.8
$ws.Cells.Item(306, 15).Value2 = 4
$ws.Cells.Item(306, 18).Value2 = 2.025
$ws.Cells.Item(306, 19).Value2 = 1.825
$ws.Cells.Item(306, 20).Value2 = 3
$ws.Cells.Item(306, 21).Value2 = 1.975
$ws.Cells.Item(306, 22).Value2 = 1.875
$ws.Cells.Item(306, 23).Value2 = 0.8
$ws.Cells.Item(306, 26).Value2 = 0.5125
$ws.Cells.Item(306, 28).Value2 = 0
$ws.Cells.Item(306, 29).Value2 = -0

# Row 307
$ws.Cells.Item(307, 2).Value2 = 4811402
$ws.Cells.Item(307, 6).Value2 = "Midtjylland"
$ws.Cells.Item(307, 7).Value2 = "Randers FC"
$ws.Cells.Item(307, 9).Value2 = 2
$ws.Cells.Item(307, 11).Value2 = 1.7
$ws.Cells.Item(307, 12).Value2 = 3.8
$ws.Cells.Item(307, 13).Value2 = 5
$ws.Cells.Item(307, 14).Value2 = 1.727
$ws.Cells.Item(307, 15).Value2 = 4.5
$ws.Cells.Item(307, 16).Value2 = 4.2
$ws.Cells.Item(307, 17).Value2 = -0.75
$ws.Cells.Item(307, 18).Value2 = 1.9
$ws.Cells.Item(307, 19).Value2 = 1.95
$ws.Cells.Item(307, 20).Value2 = 3.25
$ws.Cells.Item(307, 21).Value2 = 1.9
$ws.Cells.Item(307, 22).Value2 = 1.95
$ws.Cells.Item(307, 23).Value2 = 0.7270000000000001
$ws.Cells.Item(307, 26).Value2 = 0.45
$ws.Cells.Item(307, 27).Value2 = -0.5
$ws.Cells.Item(307, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(307, 29).Value2 = -1

# Row 352
$ws.Cells.Item(352, 2).Value2 = 5143714
$ws.Cells.Item(352, 6).Value2 = "Lyngby"
$ws.Cells.Item(352, 7).Value2 = "Randers FC"
$ws.Cells.Item(352, 8).Value2 = 0
$ws.Cells.Item(352, 11).Value2 = 3.2
$ws.Cells.Item(352, 13).Value2 = 2.1
$ws.Cells.Item(352, 14).Value2 = 2.625
$ws.Cells.Item(352, 16).Value2 = 2.625
$ws.Cells.Item(352, 17).Value2 = 0
$ws.Cells.Item(352, 18).Value2 = 1.875
$ws.Cells.Item(352, 19).Value2 = 1.975
$ws.Cells.Item(352, 20).Value2 = 2.75
$ws.Cells.Item(352, 21).Value2 = 1.9
$ws.Cells.Item(352, 22).Value2 = 1.95
$ws.Cells.Item(352, 25).Value2 = 1.625
$ws.Cells.Item(352, 27).Value2 = 0.9750000000000001
$ws.Cells.Item(352, 28).Value2 = -1
$ws.Cells.Item(352, 29).Value2 = 0.95

# Row 353
$ws.Cells.Item(353, 2).Value2 = 5143657
$ws.Cells.Item(353, 6).Value2 = "Odense BK"
$ws.Cells.Item(353, 7).Value2 = "Viborg"
$ws.Cells.Item(353, 8).Value2 = 1
$ws.Cells.Item(353, 11).Value2 = 2.5
$ws.Cells.Item(353, 13).Value2 = 2.5
$ws.Cells.Item(353, 14).Value2 = 3
$ws.Cells.Item(353, 16).Value2 = 2.25
$ws.Cells.Item(353, 17).Value2 = 0.25
$ws.Cells.Item(353, 18).Value2 = 1.86
$ws.Cells.Item(353, 19).Value2 = 2.04
$ws.Cells.Item(353, 20).Value2 = 2.5
$ws.Cells.Item(353, 21).Value2 = 1.85
$ws.Cells.Item(353, 22).Value2 = 2
$ws.Cells.Item(353, 25).Value2 = 1.25
$ws.Cells.Item(353, 27).Value2 = 1.04
$ws.Cells.Item(353, 28).Value2 = 0.8500000000000001
$ws.Cells.Item(353, 29).Value2 = -1

# Row 448
$ws.Cells.Item(448, 2).Value2 = 6461434
$ws.Cells.Item(448, 6).Value2 = "Midtjylland"
$ws.Cells.Item(448, 7).Value2 = "Lyngby"
$ws.Cells.Item(448, 10).Value2 = "H"
$ws.Cells.Item(448, 8).Value2 = 1
$ws.Cells.Item(448, 9).Value2 = 0
$ws.Cells.Item(448, 11).Value2 = 1.571
$ws.Cells.Item(448, 12).Value2 = 3.8
$ws.Cells.Item(448, 13).Value2 = 5
$ws.Cells.Item(448, 14).Value2 = 1.571
$ws.Cells.Item(448, 15).Value2 = 4
$ws.Cells.Item(448, 16).Value2 = 5.75
$ws.Cells.Item(448, 17).Value2 = -1
$ws.Cells.Item(448, 18).Value2 = 1.975
$ws.Cells.Item(448, 19).Value2 = 1.875
$ws.Cells.Item(448, 21).Value2 = 1.925
$ws.Cells.Item(448, 22).Value2 = 1.925
$ws.Cells.Item(448, 23).Value2 = 0.571
$ws.Cells.Item(448, 25).Value2 = -1
$ws.Cells.Item(448, 26).Value2 = 0
$ws.Cells.Item(448, 27).Value2 = -0
$ws.Cells.Item(448, 28).Value2 = -1
$ws.Cells.Item(448, 29).Value2 = 0.925

# Row 449
$ws.Cells.Item(449, 2).Value2 = 6461433
$ws.Cells.Item(449, 6).Value2 = "AaB"
$ws.Cells.Item(449, 7).Value2 = "Odense BK"
$ws.Cells.Item(449, 10).Value2 = "A"
$ws.Cells.Item(449, 8).Value2 = 2
$ws.Cells.Item(449, 9).Value2 = 3
$ws.Cells.Item(449, 11).Value2 = 2.375
$ws.Cells.Item(449, 12).Value2 = 3.4
$ws.Cells.Item(449, 13).Value2 = 2.8
$ws.Cells.Item(449, 14).Value2 = 1.85
$ws.Cells.Item(449, 15).Value2 = 3.75
$ws.Cells.Item(449, 16).Value2 = 4.2
$ws.Cells.Item(449, 17).Value2 = -0.5
$ws.Cells.Item(449, 18).Value2 = 1.85
$ws.Cells.Item(449, 19).Value2 = 2
$ws.Cells.Item(449, 21).Value2 = 1.975
$ws.Cells.Item(449, 22).Value2 = 1.875
$ws.Cells.Item(449, 23).Value2 = -1
$ws.Cells.Item(449, 25).Value2 = 3.2
$ws.Cells.Item(449, 26).Value2 = -1
$ws.Cells.Item(449, 27).Value2 = 1
$ws.Cells.Item(449, 28).Value2 = 0.9750000000000001
$ws.Cells.Item(449, 29).Value2 = -1

# Row 454
$ws.Cells.Item(454, 2).Value2 = 6445249
$ws.Cells.Item(454, 6).Value2 = "Odense BK"
$ws.Cells.Item(454, 7).Value2 = "Silkeborg IF"
$ws.Cells.Item(454, 9).Value2 = 0
$ws.Cells.Item(454, 11).Value2 = 2.75
$ws.Cells.Item(454, 12).Value2 = 3.5
$ws.Cells.Item(454, 13).Value2 = 2.45
$ws.Cells.Item(454, 14).Value2 = 2.45
$ws.Cells.Item(454, 15).Value2 = 3.6
$ws.Cells.Item(454, 16).Value2 = 2.7
$ws.Cells.Item(454, 17).Value2 = 0
$ws.Cells.Item(454, 18).Value2 = 1.825
$ws.Cells.Item(454, 19).Value2 = 2.025
$ws.Cells.Item(454, 20).Value2 = 2.75
$ws.Cells.Item(454, 21).Value2 = 1.85
$ws.Cells.Item(454, 22).Value2 = 2
$ws.Cells.Item(454, 23).Value2 = 1.45
$ws.Cells.Item(454, 26).Value2 = 0.825
$ws.Cells.Item(454, 28).Value2 = -1
$ws.Cells.Item(454, 29).Value2 = 1

# Row 455
$ws.Cells.Item(455, 2).Value2 = 6478386
$ws.Cells.Item(455, 6).Value2 = "Lyngby"
$ws.Cells.Item(455, 7).Value2 = "AC Horsens"
$ws.Cells.Item(455, 9).Value2 = 1
$ws.Cells.Item(455, 11).Value2 = 2.05
$ws.Cells.Item(455, 12).Value2 = 3.4
$ws.Cells.Item(455, 13).Value2 = 3.5
$ws.Cells.Item(455, 14).Value2 = 2.15
$ws.Cells.Item(455, 15).Value2 = 3.4
$ws.Cells.Item(455, 16).Value2 = 3.4
$ws.Cells.Item(455, 17).Value2 = -0.25
$ws.Cells.Item(455, 18).Value2 = 1.875
$ws.Cells.Item(455, 19).Value2 = 1.975
$ws.Cells.Item(455, 20).Value2 = 2.25
$ws.Cells.Item(455, 21).Value2 = 1.875
$ws.Cells.Item(455, 22).Value2 = 1.975
$ws.Cells.Item(455, 23).Value2 = 1.15
$ws.Cells.Item(455, 26).Value2 = 0.875
$ws.Cells.Item(455, 28).Value2 = 0.875
$ws.Cells.Item(455, 29).Value2 = -1

# Row 490
$ws.Cells.Item(490, 2).Value2 = 6478389
$ws.Cells.Item(490, 6).Value2 = "Lyngby"
$ws.Cells.Item(490, 7).Value2 = "AaB"
$ws.Cells.Item(490, 10).Value2 = "H"
$ws.Cells.Item(490, 8).Value2 = 2
$ws.Cells.Item(490, 9).Value2 = 1
$ws.Cells.Item(490, 11).Value2 = 2.75
$ws.Cells.Item(490, 12).Value2 = 3.7
$ws.Cells.Item(490, 13).Value2 = 2.3
$ws.Cells.Item(490, 14).Value2 = 3.2
$ws.Cells.Item(490, 16).Value2 = 2.1
$ws.Cells.Item(490, 17).Value2 = 0.25
$ws.Cells.Item(490, 18).Value2 = 2.025
$ws.Cells.Item(490, 19).Value2 = 1.825
$ws.Cells.Item(490, 20).Value2 = 2.5
$ws.Cells.Item(490, 21).Value2 = 1.875
$ws.Cells.Item(490, 22).Value2 = 1.975
$ws.Cells.Item(490, 23).Value2 = 2.2
$ws.Cells.Item(490, 24).Value2 = -1
$ws.Cells.Item(490, 26).Value2 = 1.025
$ws.Cells.Item(490, 28).Value2 = 0.875

# Row 491
$ws.Cells.Item(491, 2).Value2 = 6445255
$ws.Cells.Item(491, 6).Value2 = "Silkeborg IF"
$ws.Cells.Item(491, 7).Value2 = "Midtjylland"
$ws.Cells.Item(491, 10).Value2 = "D"
$ws.Cells.Item(491, 8).Value2 = 3
$ws.Cells.Item(491, 9).Value2 = 3
$ws.Cells.Item(491, 11).Value2 = 2.8
$ws.Cells.Item(491, 12).Value2 = 3.75
$ws.Cells.Item(491, 13).Value2 = 2.2
$ws.Cells.Item(491, 14).Value2 = 4
$ws.Cells.Item(491, 16).Value2 = 1.85
$ws.Cells.Item(491, 17).Value2 = 0.5
$ws.Cells.Item(491, 18).Value2 = 2.05
$ws.Cells.Item(491, 19).Value2 = 1.85
$ws.Cells.Item(491, 20).Value2 = 2.75
$ws.Cells.Item(491, 21).Value2 = 1.85
$ws.Cells.Item(491, 22).Value2 = 2
$ws.Cells.Item(491, 23).Value2 = -1
$ws.Cells.Item(491, 24).Value2 = 2.75
$ws.Cells.Item(491, 26).Value2 = 1.05
$ws.Cells.Item(491, 28).Value2 = 0.8500000000000001

# Row 499
$ws.Cells.Item(499, 2).Value2 = 6438715
$ws.Cells.Item(499, 6).Value2 = "AGF Aarhus"
$ws.Cells.Item(499, 7).Value2 = "Brondby"
$ws.Cells.Item(499, 8).Value2 = 3
$ws.Cells.Item(499, 9).Value2 = 3
$ws.Cells.Item(499, 11).Value2 = 2
$ws.Cells.Item(499, 12).Value2 = 3.5
$ws.Cells.Item(499, 13).Value2 = 3.25
$ws.Cells.Item(499, 14).Value2 = 1.666
$ws.Cells.Item(499, 16).Value2 = 4.75
$ws.Cells.Item(499, 17).Value2 = -0.75
$ws.Cells.Item(499, 18).Value2 = 1.85
$ws.Cells.Item(499, 19).Value2 = 2
$ws.Cells.Item(499, 20).Value2 = 2.75
$ws.Cells.Item(499, 21).Value2 = 1.9
$ws.Cells.Item(499, 22).Value2 = 1.95
$ws.Cells.Item(499, 26).Value2 = -1
$ws.Cells.Item(499, 27).Value2 = 1
$ws.Cells.Item(499, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(499, 29).Value2 = -1

# Row 500
$ws.Cells.Item(500, 2).Value2 = 6437830
$ws.Cells.Item(500, 6).Value2 = "FC Nordsjaelland"
$ws.Cells.Item(500, 7).Value2 = "Viborg"
$ws.Cells.Item(500, 8).Value2 = 0
$ws.Cells.Item(500, 9).Value2 = 0
$ws.Cells.Item(500, 11).Value2 = 2.25
$ws.Cells.Item(500, 12).Value2 = 3.75
$ws.Cells.Item(500, 13).Value2 = 2.75
$ws.Cells.Item(500, 14).Value2 = 2.05
$ws.Cells.Item(500, 16).Value2 = 3
$ws.Cells.Item(500, 17).Value2 = -0.25
$ws.Cells.Item(500, 18).Value2 = 1.875
$ws.Cells.Item(500, 19).Value2 = 1.975
$ws.Cells.Item(500, 20).Value2 = 3.25
$ws.Cells.Item(500, 21).Value2 = 2.05
$ws.Cells.Item(500, 22).Value2 = 1.8
$ws.Cells.Item(500, 26).Value2 = -0.5
$ws.Cells.Item(500, 27).Value2 = 0.4875
$ws.Cells.Item(500, 28).Value2 = -1
$ws.Cells.Item(500, 29).Value2 = 0.8

# Row 539
$ws.Cells.Item(539, 2).Value2 = 6779623
$ws.Cells.Item(539, 6).Value2 = "Silkeborg IF"
$ws.Cells.Item(539, 7).Value2 = "Hvidovre IF"
$ws.Cells.Item(539, 10).Value2 = "H"
$ws.Cells.Item(539, 9).Value2 = 0
$ws.Cells.Item(539, 11).Value2 = 1.533
$ws.Cells.Item(539, 12).Value2 = 4
$ws.Cells.Item(539, 13).Value2 = 6
$ws.Cells.Item(539, 14).Value2 = 1.444
$ws.Cells.Item(539, 15).Value2 = 4.75
$ws.Cells.Item(539, 16).Value2 = 7
$ws.Cells.Item(539, 17).Value2 = -1.25
$ws.Cells.Item(539, 18).Value2 = 1.875
$ws.Cells.Item(539, 19).Value2 = 1.975
$ws.Cells.Item(539, 20).Value2 = 3
$ws.Cells.Item(539, 21).Value2 = 1.9
$ws.Cells.Item(539, 22).Value2 = 1.95
$ws.Cells.Item(539, 23).Value2 = 0.444
$ws.Cells.Item(539, 24).Value2 = -1
$ws.Cells.Item(539, 26).Value2 = -0.5
$ws.Cells.Item(539, 27).Value2 = 0.4875
$ws.Cells.Item(539, 29).Value2 = 0.95

# Row 540
$ws.Cells.Item(540, 2).Value2 = 6779624
$ws.Cells.Item(540, 6).Value2 = "Lyngby"
$ws.Cells.Item(540, 7).Value2 = "FC Nordsjaelland"
$ws.Cells.Item(540, 10).Value2 = "D"
$ws.Cells.Item(540, 9).Value2 = 1
$ws.Cells.Item(540, 11).Value2 = 4.5
$ws.Cells.Item(540, 12).Value2 = 3.6
$ws.Cells.Item(540, 13).Value2 = 1.75
$ws.Cells.Item(540, 14).Value2 = 4.5
$ws.Cells.Item(540, 15).Value2 = 3.6
$ws.Cells.Item(540, 16).Value2 = 1.8
$ws.Cells.Item(540, 17).Value2 = 0.75
$ws.Cells.Item(540, 18).Value2 = 1.825
$ws.Cells.Item(540, 19).Value2 = 2.025
$ws.Cells.Item(540, 20).Value2 = 2.75
$ws.Cells.Item(540, 21).Value2 = 2
$ws.Cells.Item(540, 22).Value2 = 1.85
$ws.Cells.Item(540, 23).Value2 = -1
$ws.Cells.Item(540, 24).Value2 = 2.6
$ws.Cells.Item(540, 26).Value2 = 0.825
$ws.Cells.Item(540, 27).Value2 = -1
$ws.Cells.Item(540, 29).Value2 = 0.8500000000000001
